$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 (RF) updated metrics
$ws.Range("B3").Value = 0.793
$ws.Range("C3").Value = 0.886
$ws.Range("D3").Value = 0.647
$ws.Range("E3").Value = 0.901
$ws.Range("F3").Value = 0.893
$ws.Range("G3").Value = 0.108
$ws.Range("H3").Value = 0.329
$ws.Range("I3").Value = 0.238
$ws.Range("J3").Value = 0.973

# Row 4 (NN) updated metrics
$ws.Range("E4").Value = 0.636
$ws.Range("F4").Value = 0.606
$ws.Range("G4").Value = 0.398
$ws.Range("H4").Value = 0.631
$ws.Range("I4").Value = 0.477
$ws.Range("J4").Value = 0.898

# Row 5 (RNN) updated metrics
$ws.Range("E5").Value = 0.5590000000000001
$ws.Range("F5").Value = 0.54
$ws.Range("G5").Value = 0.481
$ws.Range("H5").Value = 0.694
$ws.Range("I5").Value = 0.52
$ws.Range("J5").Value = 0.836
